$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.86469999999999
$ws.Range("A9").Value = -20.22139999999997
$ws.Range("A18").Value = -23.10550000000001
$ws.Range("A20").Value = -22.25540000000001
$ws.Range("D21").Value = -7.5144
